# Generate Report for Handback
#
# The handback transform for bc562813-ddb9-4366-b380-317f741ad6da failed:
# the handback file name did not match the expected handoff-derived name.
# Update the status on the Overview sheet and on each language sheet
# (zh-cn, de-de), and record the mismatch detail in the "Error Detail"
# column (P) for each language sheet. Also widen the Error Detail column
# so the new message is readable.

$wb = $excel.ActiveWorkbook

$statusNew = "Handback transform failed"

$zhCnDetail = "Handback file name: 211ghmxa.2ki is different with handoff file name: bc562813-ddb9-4366-b380-317f741ad6da.726c1488d5d7c2a8fc4ac281ee70e16eb8c02be2.zh-cn."
$deDeDetail = "Handback file name: 211ghmxa.2ki is different with handoff file name: bc562813-ddb9-4366-b380-317f741ad6da.726c1488d5d7c2a8fc4ac281ee70e16eb8c02be2.de-de."

# --- Overview sheet: update the Status + Status(duplicate) columns for the
#     bc562813 row (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusNew
$overview.Range("F3").Value = $statusNew

# --- zh-cn sheet: update Status (col C) and Error Detail (col P) for the
#     bc562813 row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusNew
$zhcn.Range("P3").Value = $zhCnDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: update Status (col C) and Error Detail (col P) for the
#     bc562813 row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusNew
$dede.Range("P3").Value = $deDeDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
